$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the forest floor consumption model: J11 and J13 used to be the
# constant 0.5 emission-factor adjustment; the text equations are now
# driven directly off the stated 1.4 (g consumed / g available) ratio.
$ws.Range("J11").Formula = "=1/1.4"
$ws.Range("J13").Formula = "=1/1.4"

# Move the active selection to J12 (cosmetic, matches the saved view state)
$ws.Range("J12").Select()
